$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.784.97"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "2.093.69"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.46"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.15"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.37"
$ws.Range("E12").Value = "  +5.11%  "
$ws.Range("D13").Value = "2.407.64"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.09"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.806"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.51"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "2.083.97"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "38.697.47"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.95"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.90"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.87"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("E28").Value = "  +5.79%  "
$ws.Range("E29").Value = "  +6.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.55"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("E31").Value = "  +2.35%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0611"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.54"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  +1.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.14"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0229"
$ws.Range("E41").Value = "  +5.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.98"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "1.532.00"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0913"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.71"
$ws.Range("E46").Value = "  +6.58%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.14"
$ws.Range("E47").Value = "  +1.98%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D51").Value = "2.288.93"
$ws.Range("E51").Value = "  +0.08%  "
